$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ДЗ_1 (column C) and ДЗ_2 (column D) grades of 5 recorded for several students
$ws.Range("C6").Value = 5
$ws.Range("C7").Value = 5
$ws.Range("C12").Value = 5
$ws.Range("D12").Value = 5
$ws.Range("C19").Value = 5
$ws.Range("C20").Value = 5
$ws.Range("C26").Value = 5
$ws.Range("D26").Value = 5
$ws.Range("C27").Value = 5
$ws.Range("C29").Value = 5
$ws.Range("C31").Value = 5
$ws.Range("D31").Value = 5

# Move the active selection to D30, matching where the author last worked
$ws.Range("D30").Select()
